# Nexial showcase workbook update:
#   - add "aws.ses" module (sendMail / sendTextMail) to the hidden '#system'
#     sheet that feeds the macro drop-downs, and register it in the
#     defined-name table + the "target" module list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------------
# 1) Make room for the new "aws.ses" column. The module columns run B..Z in
#    alphabetical order (aws.s3, base, csv, ... xml); "aws.ses" sorts right
#    after "aws.s3", so insert a new column before the current column C and
#    everything from base onward slides one column to the right (C->D, D->E,
#    ... Z->AA).
# ---------------------------------------------------------------------------
$ws.Columns("C:C").Insert()

# ---------------------------------------------------------------------------
# 2) The "target" list (column A) is the alphabetical list of module names;
#    "aws.ses" is inserted between "aws.s3" (row 2) and "base" (row 3), so
#    push the existing rows 3..26 down to 4..27, one cell at a time (bottom
#    up so nothing gets clobbered).
# ---------------------------------------------------------------------------
for ($r = 26; $r -ge 3; $r--) {
    $ws.Range("A" + ($r + 1)).Value2 = $ws.Range("A" + $r).Value2
}

# ---------------------------------------------------------------------------
# 3) Populate the new cells.
# ---------------------------------------------------------------------------
$ws.Range("A3").Value2 = "aws.ses"

$ws.Range("C1").Value2 = "aws.ses"
$ws.Range("C2").Value2 = "sendMail(profile,to,subject,body)"
$ws.Range("C3").Value2 = "sendTextMail(profile,to,subject,body)"

# ---------------------------------------------------------------------------
# 4) Fix up the defined names: every module from "base" onward moved one
#    column to the right, "target" grew by one row, and "aws.ses" is brand
#    new. Names that already pointed at an untouched column (date, db,
#    math, mq, nextgen) are left as-is.
# ---------------------------------------------------------------------------
$afterRanges = [ordered]@{
    "base"      = "`$D`$2:`$D`$36"
    "csv"       = "`$E`$2:`$E`$5"
    "desktop"   = "`$F`$2:`$F`$92"
    "excel"     = "`$G`$2:`$G`$14"
    "external"  = "`$H`$2:`$H`$3"
    "image"     = "`$I`$2:`$I`$5"
    "io"        = "`$J`$2:`$J`$24"
    "jms"       = "`$K`$2:`$K`$4"
    "json"      = "`$L`$2:`$L`$14"
    "mail"      = "`$M`$2:`$M`$2"
    "number"    = "`$N`$2:`$N`$15"
    "pdf"       = "`$O`$2:`$O`$16"
    "rdbms"     = "`$P`$2:`$P`$7"
    "redis"     = "`$Q`$2:`$Q`$10"
    "sms"       = "`$R`$2:`$R`$2"
    "sound"     = "`$S`$2:`$S`$5"
    "ssh"       = "`$T`$2:`$T`$9"
    "step"      = "`$U`$2:`$U`$4"
    "target"    = "`$A`$2:`$A`$27"
    "web"       = "`$V`$2:`$V`$117"
    "webalert"  = "`$W`$2:`$W`$8"
    "webcookie" = "`$X`$2:`$X`$8"
    "ws"        = "`$Y`$2:`$Y`$17"
    "ws.async"  = "`$Z`$2:`$Z`$8"
    "xml"       = "`$AA`$2:`$AA`$11"
}

foreach ($nm in $afterRanges.Keys) {
    $n = $wb.Names.Item($nm)
    $n.RefersTo = "='#system'!" + $afterRanges[$nm]
}

$wb.Names.Add("aws.ses", "='#system'!`$C`$2:`$C`$3")
